$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the revised / expanded notes. Order matters for shared-string index
# assignment (matches how the original author's edit ended up ordering the
# table): the "admin abilities" cell first, then the "layout" cell, then the
# "3pm" cell.
$ws.Range("D29").Value = "Need to decide about abilities and presentation for admin --- probably same views but with some kind of ability to see ""removed content"" stats, and with edit/remove options in more places."
$ws.Range("D28").Value = "Need to: make layout for forum, thread, admin view of questionable users"
$ws.Range("D24").Value = "3pm - 4:45 Revised project plan, user stories, and project plan to reflect current status."

# Row 29's text now wraps to two lines, so its height grows to match the
# other wrapped-text rows in this column.
$ws.Rows(29).RowHeight = 30

# The old trailing note ("Need to revise intro...") in row 30 is dropped
# entirely as part of this edit, so remove that row.
$ws.Rows(30).Delete()

# Reflect the new bottom of the used range in the active selection.
$ws.Range("D30").Select()
